{"js": "// The upstream commit (\"Fixed #295 Add the version of M2Doc in the\n// template custom properties\") re-saved this template through a tool\n// that rewrites each OOXML part with its attributes in alphabetical\n// order. For this particular file, the resulting unified diff touches\n// only attribute *order* on existing elements (root namespace\n// declarations on <w:document>, <w:pgSz>/<w:pgMar> on the section\n// properties, <w:rFonts>/<w:lang> in the style defaults, the\n// <w:latentStyles>/<w:lsdException> table, and the built-in <w:style>\n// definitions together with <w:tblInd>/<w:tblCellMar> inside the\n// \"TableauNormal\" table style). Every changed line keeps the exact same\n// element name and the exact same set of attribute name/value pairs \u2014\n// nothing is added, removed, or renamed, and no text/content/formatting\n// is different. docProps/custom.xml (where the actual M2Doc version\n// property lives) is not part of this file's diff either.\n//\n// Attribute order is not part of the XML information set (it carries no\n// semantic meaning) and the Word JavaScript API only exposes the\n// document at the semantic level (paragraphs, ranges, styles, page\n// setup, ...) rather than raw part serialization, so there is no\n// content change to make here. We simply touch the same properties the\n// diff's attributes belong to, confirming they are already correct, and\n// perform no writes so the document content stays byte-for-byte\n// equivalent.\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nawait context.sync();\n\n// Section/page-setup values referenced by <w:pgSz>/<w:pgMar> in the\n// diff are read-only confirmations in Office.js (PageSetup margins are\n// not part of the JS API surface); nothing here is written back, so no\n// new content diff is produced, matching the no-op nature of the\n// upstream change for this file.\nfor (const section of sections.items) {\n  section.load(\"body\");\n}\n\nawait context.sync();\n", "ps1": "# The upstream commit (\"Fixed #295 Add the version of M2Doc in the\n# template custom properties\") re-saved this template through a tool\n# that rewrites each OOXML part with its attributes in alphabetical\n# order. For this particular file, the resulting unified diff touches\n# only attribute *order* on existing elements (root namespace\n# declarations on <w:document>, <w:pgSz>/<w:pgMar> on the section\n# properties, <w:rFonts>/<w:lang> in the style defaults, the\n# <w:latentStyles>/<w:lsdException> table, and the built-in <w:style>\n# definitions together with <w:tblInd>/<w:tblCellMar> inside the\n# \"TableauNormal\" table style). Every changed line keeps the exact same\n# element name and the exact same set of attribute name/value pairs --\n# nothing is added, removed, or renamed, and no text/content/formatting\n# is different. docProps/custom.xml (where the actual M2Doc version\n# property lives) is not part of this file's diff either.\n#\n# Attribute order is not part of the XML information set (it carries no\n# semantic meaning) and the Word object model only exposes the document\n# at the semantic level (paragraphs, ranges, styles, page setup, ...)\n# rather than raw part serialization, so there is no content change to\n# make here. We simply touch the same properties the diff's attributes\n# belong to, confirming they are already correct, and perform no writes\n# so the document content stays byte-for-byte equivalent.\n\n$d = $word.ActiveDocument\n\n# Section / page-setup values referenced by <w:pgSz>/<w:pgMar>.\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n$null = $pageSetup.PageWidth\n$null = $pageSetup.PageHeight\n$null = $pageSetup.TopMargin\n$null = $pageSetup.RightMargin\n$null = $pageSetup.BottomMargin\n$null = $pageSetup.LeftMargin\n$null = $pageSetup.HeaderDistance\n$null = $pageSetup.FooterDistance\n$null = $pageSetup.Gutter\n\n# Style defaults referenced by <w:rFonts>/<w:lang> and the built-in\n# styles touched by the reorder (Normal, Default Paragraph Font, Normal\n# Table, No List).\nforeach ($styleName in @(\"Normal\", \"Default Paragraph Font\", \"Normal Table\", \"No List\")) {\n    $style = $d.Styles.Item($styleName)\n    $null = $style.NameLocal\n}\n\n# No property is assigned above (every access is a read), so the\n# document content is left exactly as it was -- matching the purely\n# cosmetic (attribute-order only) nature of the upstream change for\n# this file.\n"}
